{"js": "// Module 11 case study edit \u2014 three separate text rewrites inside the\n// \"Unit Assessment\" document, located by searching for the unique\n// sentence fragments shown in the diff and replacing them with the\n// post-edit wording (the content is identical once the runs that the\n// diff spread the new wording across are concatenated back together).\n\nconst body = context.document.body;\n\n// --- Edit 1 -------------------------------------------------------\n// \"...is diagnosed, these therapies have to be \" ->\n// \"...is diagnosed with Parkinson\u2019s disease, these therapies are\n//  critical, and have to be \"\nconst search1 = body.search(\n  \"Knowing that 70% of these neurons are gone by the time a patient is diagnosed, these therapies have to be \",\n  { matchCase: true }\n);\nsearch1.load(\"items\");\nawait context.sync();\n\nif (search1.items.length === 0) {\n  throw new Error(\"Edit 1: target text not found\");\n}\nsearch1.items[0].insertText(\n  \"Knowing that 70% of these neurons are gone by the time a patient is diagnosed\" +\n    \" with Parkinson\\u2019s disease, these therapies are critical\" +\n    \", \" +\n    \"and\" +\n    \" have to be \",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- Edit 2 -------------------------------------------------------\n// \"...expose the mice to a battery of tests since...\" ->\n// \"...expose the mice to various types of tests since...\"\nconst search2 = body.search(\n  \"It is important to expose the mice to a battery of tests since it has been observed that mice can display impairments on a specific test but \",\n  { matchCase: true }\n);\nsearch2.load(\"items\");\nawait context.sync();\n\nif (search2.items.length === 0) {\n  throw new Error(\"Edit 2: target text not found\");\n}\nsearch2.items[0].insertText(\n  \"It is important to expose the mice to \" +\n    \"various types of \" +\n    \"tests since it has been observed that mice can display impairments on a specific test but \",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- Edit 3 -------------------------------------------------------\n// \"...decrease uncertainty ... with PD.  More rigorous tests are needed\n//  to evaluate...\" ->\n// \"...decrease the uncertainty ... with PD and be declared an efficient\n//  therapy.  More rigorous tests on animal models are required to\n//  evaluate...\"\nconst search3 = body.search(\n  \" uncertainty when contemplating applying this strategy to patients with PD.  More rigorous tests are needed to evaluate long-term impact of reprogrammed neurons on the brain and overall patient condition.\",\n  { matchCase: true }\n);\nsearch3.load(\"items\");\nawait context.sync();\n\nif (search3.items.length === 0) {\n  throw new Error(\"Edit 3: target text not found\");\n}\nsearch3.items[0].insertText(\n  \" \" +\n    \"the \" +\n    \"uncertainty when contemplating applying this strategy to patients with PD\" +\n    \" and be declared an efficient therapy\" +\n    \".  More rigorous tests \" +\n    \"o\" +\n    \"n animal models \" +\n    \"are \" +\n    \"req\" +\n    \"u\" +\n    \"ired\" +\n    \" to evaluate long-term impact of reprogrammed neurons on the brain and overall patient condition.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Module 11 case study edit \u2014 three separate text rewrites inside the\n# \"Unit Assessment\" document. Each unique sentence fragment is located\n# with Find/Replace and rewritten to the post-edit wording (the content\n# is identical once the runs the diff spread the new wording across are\n# concatenated back together).\n\n$d = $word.ActiveDocument\n\n# --- Edit 1 ---------------------------------------------------------\n# \"...is diagnosed, these therapies have to be \" ->\n# \"...is diagnosed with Parkinson\u2019s disease, these therapies are\n#  critical, and have to be \"\n$find1 = \"Knowing that 70% of these neurons are gone by the time a patient is diagnosed, these therapies have to be \"\n$replace1 = \"Knowing that 70% of these neurons are gone by the time a patient is diagnosed with Parkinson\u2019s disease, these therapies are critical, and have to be \"\n\n$r1 = $d.Content\n$r1.Find.ClearFormatting()\n$r1.Find.Execute($find1, $false, $false, $false, $false, $false, $true, 1, $false, $replace1, 2)\n\n# --- Edit 2 ---------------------------------------------------------\n# \"...expose the mice to a battery of tests since...\" ->\n# \"...expose the mice to various types of tests since...\"\n$find2 = \"It is important to expose the mice to a battery of tests since it has been observed that mice can display impairments on a specific test but \"\n$replace2 = \"It is important to expose the mice to various types of tests since it has been observed that mice can display impairments on a specific test but \"\n\n$r2 = $d.Content\n$r2.Find.ClearFormatting()\n$r2.Find.Execute($find2, $false, $false, $false, $false, $false, $true, 1, $false, $replace2, 2)\n\n# --- Edit 3 ---------------------------------------------------------\n# \"...decrease uncertainty ... with PD.  More rigorous tests are needed\n#  to evaluate...\" ->\n# \"...decrease the uncertainty ... with PD and be declared an efficient\n#  therapy.  More rigorous tests on animal models are required to\n#  evaluate...\"\n$find3 = \" uncertainty when contemplating applying this strategy to patients with PD.  More rigorous tests are needed to evaluate long-term impact of reprogrammed neurons on the brain and overall patient condition.\"\n$replace3 = \" the uncertainty when contemplating applying this strategy to patients with PD and be declared an efficient therapy.  More rigorous tests on animal models are required to evaluate long-term impact of reprogrammed neurons on the brain and overall patient condition.\"\n\n$r3 = $d.Content\n$r3.Find.ClearFormatting()\n$r3.Find.Execute($find3, $false, $false, $false, $false, $false, $true, 1, $false, $replace3, 2)\n"}
